$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet previously held data in columns A:I (rows 1-4). Columns J:N
# are new - they duplicate the existing E:I columns for each of the three
# data rows (2-4), as if the last five data columns were copied and
# pasted one block to the right.

# Row 2 ("sprint")
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 8

# Row 3 ("Duration")
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 17
$ws.Range("L3").Value = 15
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 15

# Row 4 ("Grooming")
$ws.Range("J4").Value = 0.5
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.5
$ws.Range("N4").Value = 0.5

# Copy the formatting (borders/fill/alignment) of the source block E2:I4
# onto the newly filled J2:N4 block, matching a copy/paste of that range.
$ws.Range("E2:I4").Copy() | Out-Null
$null = $ws.Range("J2").PasteSpecial(-4122)

# Excel leaves the pasted range selected after a paste operation.
$null = $ws.Range("J2:N4").Select()
